$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellStyle($addr, $styleAddr) {
    $ws.Range($styleAddr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# --- Row 68 ---
$ws.Range("Z68").Value = 'El usuario valida que se muestre la pestaña "Movimientos".||El sistema muestra la sección Movimientos.'
$ws.Range("AA68").Value = 'El usuario valida que se muestre el campo de búsqueda "Seleccione fecha inicio".||El sistema muestra el campo de fecha habilitado.'
$ws.Range("AB68").Value = 'El usuario valida que se muestre el campo de búsqueda "Seleccione fecha fin".||El sistema muestra el campo de fecha deshabilitado.'
$ws.Range("AC68").Value = 'El usuario valida que se muestre el botón de Descarga.||El sistema muestra el botón de "Descarga" de forma correcta.'
$ws.Range("AD68").Value = 'El usuario valida que se muestre la columna "Fecha".||El sistema muestra la columna "Fecha" en formato DD/MMM/AAAA".'
$ws.Range("AE68").Value = 'El usuario valida que se muestre la columna "Número de autorización".||El sistema muestra la columna "Número de autorización" de forma correcta.'
$ws.Range("AF68").Value = 'El usuario valida que se muestre la columna "Transacción".||El sistema muestra la columna "Transacción" de forma correcta.'
$ws.Range("AG68").Value = 'El usuario valida que se muestre la columna "Monto".||El sistema muestra la columna "Monto" del movimiento con el siguiente formato:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales
Se muestra el valor positivo si se trata de un Abono.
Se muestra el valor negativo si se trata de un Cargo.'
$ws.Range("AH68").Value = 'El usuario valida que se muestre la columna "Tipo de operación".||El sistema muestra la columna "Tipo de operación" con los valores "Cargo" o "Abono".'
$ws.Range("AI68").Value = 'El usuario valida que se muestre la columna opciones (…).||El sistema la columna opciones (…) de forma correcta.'
$ws.Range("AJ68").Value = 'El usuario valida que se muestre la pestaña "Movimientos".||El sistema muestra la sección Movimientos.'
$ws.Range("AK68").Value = 'El usuario indica una fecha de inicio en el campo de búsqueda "Seleccione fecha inicio".||El sistema muestra la fecha de inicio de forma correcta.'
$ws.Range("AL68").Value = 'El usuario indica una fecha fin en el campo de búsqueda "Seleccione fecha fin".||El sistema muestra la fecha fin de forma correcta y el resultado de la búsqueda de movimientos en la siguiente tabla:
Columna "Fecha"
Columna "Número de autorización"
Columna "Transacción"
Columna "Monto"
Columna "Tipo de operación"'
$ws.Range("AM68").Value = 'El usuario da clic en el botón "Opciones" con el icono (…).||El sistema muestra la pantalla modal de forma correcta.'
$ws.Range("AN68").Value = 'El usuario da clic en el botón "Imprimir".||El sistema muestra la descarga del archivo PDF de forma correcta.'

Set-CellStyle "Z68" "H10"
Set-CellStyle "AA68" "H10"
Set-CellStyle "AB68" "H10"
Set-CellStyle "AC68" "O15"
Set-CellStyle "AD68" "O15"
Set-CellStyle "AE68" "A2"
Set-CellStyle "AF68" "O15"
Set-CellStyle "AG68" "O15"
Set-CellStyle "AH68" "A2"
Set-CellStyle "AI68" "O15"
Set-CellStyle "AJ68" "A2"
Set-CellStyle "AK68" "O15"
Set-CellStyle "AL68" "A2"
Set-CellStyle "AM68" "A2"
Set-CellStyle "AN68" "A2"

# --- Row 69 ---
$ws.Range("A69").Value = 'TC_002_68_administradorConsulta'
$ws.Range("C69").Value = 'El usuario despliega la página de CCOP mediante el siguiente URL:
https://ccop-u.scointnet.net?countryId=MX||El sistema muestra la página de CCOP de forma correcta.'
$ws.Range("D69").Value = 'El usuario ingresa el siguiente dato: Nombre de usuario||El sistema permite ingresar el nombre de usuario de forma correcta.'
$ws.Range("E69").Value = 'El usuario ingresa el siguiente dato: Contraseña||El sistema permite ingresar la contraseña de forma correcta.'
$ws.Range("F69").Value = 'El usuario ingresa el siguiente dato: Token||El sistema permite ingresar el token de forma correcta.'
$ws.Range("G69").Value = 'El usuario da clic en el botón "Ingresar"||El sistema muestra la página principal de CCOP "Resumen consolidado de productos".'
$ws.Range("H69").Value = 'El usuario da clic en la pestaña "Productos y servicios"||El sistema muestra la página "Cuentas".'
$ws.Range("I69").Value = 'El usuario da clic en la pestaña "Cuentas Cheque USD"||El sistema muestra la lista de cuentas de cheques USD.'
$ws.Range("J69").Value = 'El usuario da clic en una cuenta cheque USD.||El sistema muestra el detalle de la cuenta cheque USD.'
$ws.Range("K69").Value = 'El usuario valida el campo "Número de cuenta cheque USD".||El sistema muestra el número de la cuenta cheque USD de forma correcta.'
$ws.Range("L69").Value = 'El usuario valida el campo "Identificación cliente".||El sistema muestra la identificación del cliente de forma correcta.'
$ws.Range("M69").Value = 'El usuario valida el campo "Ciudad".||El sistema muestra la ciudad de forma correcta.'
$ws.Range("N69").Value = 'El usuario valida el Nombre de cliente.||El sistema muestra el nombre del cliente de forma correcta.'
$ws.Range("O69").Value = 'El usuario valida el estado de la cuenta.||El sistema muestra el estado de la cuenta en la parte superior derecha como "ACTIVA".'
$ws.Range("P69").Value = 'El usuario valida que se muestre la columna "Saldo total".||El sistema muestra la columna "Saldo total" de forma correcta con formato de moneda:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("Q69").Value = 'El usuario valida que se muestre la columna "Saldo disponible"||El sistema muestra la columna "Saldo disponible" de forma correcta con formato de moneda:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("R69").Value = 'El usuario valida que se muestre la columna "Saldo salvo buen cobro"||El sistema muestra la columna "Saldo salvo buen cobro" de forma correcta con formato de moneda:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("S69").Value = 'El usuario valida que se muestre la columna "Saldo límite de sobregiro"||El sistema muestra la columna "Saldo límite de sobregiro" de forma correcta con formato de moneda:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("T69").Value = 'El usuario valida que se muestre el campo "Alias".||El sistema muestra el campo Alias de forma correcta.'
$ws.Range("U69").Value = 'El usuario da clic en la opción "Ocultar detalles".||El sistema permite que la opción "Ocultar detalles" oculte la sección de Saldos.'
$ws.Range("V69").Value = 'El usuario valida que se muestre la pestaña "Movimientos".||El sistema muestra la sección Movimientos.'
$ws.Range("W69").Value = 'El usuario valida que se muestre el campo de búsqueda "Seleccione fecha inicio".||El sistema muestra el campo de fecha habilitado.'
$ws.Range("X69").Value = 'El usuario valida que se muestre el campo de búsqueda "Seleccione fecha fin".||El sistema muestra el campo de fecha deshabilitado.'
$ws.Range("Y69").Value = 'El usuario valida que se muestre el botón "Descargar reporte".||El sistema muestra el botón de "Descargar reporte" de forma correcta.'
$ws.Range("Z69").Value = 'El usuario valida que se muestre la columna "Fecha del movimiento".||El sistema muestra la columna "Fecha del movimiento" en formato DD/MMM/AAAA".'
$ws.Range("AA69").Value = 'El usuario valida que se muestre la columna "Tipo de operación".||El sistema muestra la columna "Tipo de operación" con los valores "Cargo" o "Abono".'
$ws.Range("AB69").Value = 'El usuario valida que se muestre la columna "Descripción".||El sistema muestra la columna "Descripción" con la descripción del movimiento.'
$ws.Range("AC69").Value = 'El usuario valida que se muestre la columna "Monto".||El sistema muestra la columna "Monto" del movimiento con el siguiente formato:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales
Se muestra el valor positivo si se trata de un Abono.
Se muestra el valor negativo si se trata de un Cargo.'
$ws.Range("AD69").Value = 'El usuario valida que se muestre la columna "Saldo final".||El sistema muestra la columna "Monto" del movimiento con el siguiente formato:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("AE69").Value = 'El usuario valida que se muestre la columna opciones (…).||El sistema la columna opciones (…) de forma correcta.'
$ws.Range("AF69").Value = 'El usuario valida que se muestre la pestaña "Movimientos".||El sistema muestra la sección Movimientos.'
$ws.Range("AG69").Value = 'El usuario indica una fecha de inicio en el campo de búsqueda "Seleccione fecha inicio".||El sistema muestra la fecha de inicio de forma correcta.'
$ws.Range("AH69").Value = 'El usuario indica una fecha fin en el campo de búsqueda "Seleccione fecha fin".||El sistema muestra la fecha fin de forma correcta y el resultado de la búsqueda de movimientos en la siguiente tabla:
Columna "Fecha del movimiento"
Columna "Tipo de operación"
Columna "Descripción"
Columna "Monto"
Columna "Saldo final"'
$ws.Range("AI69").Value = 'El usuario da clic en el botón "Opciones" con el icono (…).||El sistema muestra la pantalla modal de forma correcta.'
$ws.Range("AJ69").Value = 'El usuario da clic en el botón "Imprimir".||El sistema muestra la descarga del archivo PDF de forma correcta.'

Set-CellStyle "A69" "A66"
Set-CellStyle "C69" "C2"
Set-CellStyle "D69" "D2"
Set-CellStyle "E69" "D2"
Set-CellStyle "F69" "F42"
Set-CellStyle "G69" "H38"
Set-CellStyle "H69" "H51"
Set-CellStyle "I69" "A2"
Set-CellStyle "J69" "A2"
Set-CellStyle "K69" "A2"
Set-CellStyle "L69" "A2"
Set-CellStyle "M69" "H10"
Set-CellStyle "N69" "H10"
Set-CellStyle "O69" "H10"
Set-CellStyle "P69" "H10"
Set-CellStyle "Q69" "H10"
Set-CellStyle "R69" "H10"
Set-CellStyle "S69" "H10"
Set-CellStyle "T69" "H10"
Set-CellStyle "U69" "H10"
Set-CellStyle "V69" "H10"
Set-CellStyle "W69" "U42"
Set-CellStyle "X69" "H10"
Set-CellStyle "Y69" "U42"
Set-CellStyle "Z69" "H10"
Set-CellStyle "AA69" "H10"
Set-CellStyle "AB69" "H10"
Set-CellStyle "AC69" "A2"
Set-CellStyle "AD69" "O15"
Set-CellStyle "AE69" "A2"
Set-CellStyle "AF69" "O15"
Set-CellStyle "AG69" "A2"
Set-CellStyle "AH69" "A2"
Set-CellStyle "AI69" "O15"
Set-CellStyle "AJ69" "O15"

# --- Row 70 ---
$ws.Range("A70").Value = 'TC_002_69_administradorConsulta'
$ws.Range("C70").Value = 'El usuario despliega la página de CCOP mediante el siguiente URL:
https://ccop-u.scointnet.net?countryId=MX||El sistema muestra la página de CCOP de forma correcta.'
$ws.Range("D70").Value = 'El usuario ingresa el siguiente dato: Nombre de usuario||El sistema permite ingresar el nombre de usuario de forma correcta.'
$ws.Range("E70").Value = 'El usuario ingresa el siguiente dato: Contraseña||El sistema permite ingresar la contraseña de forma correcta.'
$ws.Range("F70").Value = 'El usuario ingresa el siguiente dato: Token||El sistema permite ingresar el token de forma correcta.'
$ws.Range("G70").Value = 'El usuario da clic en el botón "Ingresar"||El sistema muestra la página principal de CCOP "Resumen consolidado de productos".'
$ws.Range("H70").Value = 'El usuario da clic en la pestaña "Productos y servicios"||El sistema muestra la página "Cuentas".'
$ws.Range("I70").Value = 'El usuario valida las cuentas de cheques de la pestaña "Cuentas Cheque MXN"||El sistema muestra la lista de cuentas de cheques MXN.'
$ws.Range("J70").Value = 'El usuario da clic en una cuenta cheque MXN.||El sistema muestra el detalle de la cuenta cheque MXN.'
$ws.Range("K70").Value = 'El usuario valida el campo "Número de cuenta cheque MXN".||El sistema muestra el número de la cuenta cheque MXN de forma correcta.'
$ws.Range("L70").Value = 'El usuario valida el campo "Identificación cliente".||El sistema muestra la identificación del cliente de forma correcta.'
$ws.Range("M70").Value = 'El usuario valida el campo "Ciudad".||El sistema muestra la ciudad de forma correcta.'
$ws.Range("N70").Value = 'El usuario valida el Nombre de cliente.||El sistema muestra el nombre del cliente de forma correcta.'
$ws.Range("O70").Value = 'El usuario valida el estado de la cuenta.||El sistema muestra el estado de la cuenta en la parte superior derecha como "ACTIVA".'
$ws.Range("P70").Value = 'El usuario valida que se muestre la columna "Saldo total".||El sistema muestra la columna "Saldo total" de forma correcta con formato de moneda:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("Q70").Value = 'El usuario valida que se muestre la columna "Saldo disponible"||El sistema muestra la columna "Saldo disponible" de forma correcta con formato de moneda:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("R70").Value = 'El usuario valida que se muestre la columna "Saldo salvo buen cobro"||El sistema muestra la columna "Saldo salvo buen cobro" de forma correcta con formato de moneda:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("S70").Value = 'El usuario valida que se muestre la columna "Saldo límite de sobregiro"||El sistema muestra la columna "Saldo límite de sobregiro" de forma correcta con formato de moneda:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("T70").Value = 'El usuario valida que se muestre el campo "Alias".||El sistema muestra el campo Alias de forma correcta.'
$ws.Range("U70").Value = 'El usuario da clic en la opción "Ocultar detalles".||El sistema permite que la opción "Ocultar detalles" oculte la sección de Saldos.'
$ws.Range("V70").Value = 'El usuario valida que se muestre la pestaña "Movimientos".||El sistema muestra la sección Movimientos.'
$ws.Range("W70").Value = 'El usuario valida que se muestre el campo de búsqueda "Seleccione fecha inicio".||El sistema muestra el campo de fecha habilitado.'
$ws.Range("X70").Value = 'El usuario valida que se muestre el campo de búsqueda "Seleccione fecha fin".||El sistema muestra el campo de fecha deshabilitado.'
$ws.Range("Y70").Value = 'El usuario valida que se muestre el botón de Descarga.||El sistema muestra el botón de "Descarga" de forma correcta.'
$ws.Range("Z70").Value = 'El usuario valida que se muestre la columna "Fecha del movimiento".||El sistema muestra la columna "Fecha del movimiento" en formato DD/MMM/AAAA".'
$ws.Range("AA70").Value = 'El usuario valida que se muestre la columna "Tipo de operación".||El sistema muestra la columna "Tipo de operación" con los valores "Cargo" o "Abono".'
$ws.Range("AB70").Value = 'El usuario valida que se muestre la columna "Descripción".||El sistema muestra la columna "Descripción" con la descripción del movimiento.'
$ws.Range("AC70").Value = 'El usuario valida que se muestre la columna "Monto".||El sistema muestra la columna "Monto" del movimiento con el siguiente formato:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales
Se muestra el valor positivo si se trata de un Abono.
Se muestra el valor negativo si se trata de un Cargo.'
$ws.Range("AD70").Value = 'El usuario valida que se muestre la columna "Saldo final".||El sistema muestra la columna "Monto" del movimiento con el siguiente formato:
Signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("AE70").Value = 'El usuario valida que se muestre la columna opciones (…).||El sistema la columna opciones (…) de forma correcta.'
$ws.Range("AF70").Value = 'El usuario valida que se muestre la pestaña "Movimientos".||El sistema muestra la sección Movimientos.'
$ws.Range("AG70").Value = 'El usuario indica una fecha de inicio en el campo de búsqueda "Seleccione fecha inicio".||El sistema muestra la fecha de inicio de forma correcta.'
$ws.Range("AH70").Value = 'El usuario indica una fecha fin en el campo de búsqueda "Seleccione fecha fin".||El sistema muestra la fecha fin de forma correcta y el resultado de la búsqueda de movimientos en la siguiente tabla:
Columna "Fecha del movimiento"
Columna "Tipo de operación"
Columna "Descripción"
Columna "Monto"
Columna "Saldo final"'
$ws.Range("AI70").Value = 'El usuario da clic en el botón "Opciones" con el icono (…).||El sistema muestra la pantalla modal de forma correcta.'
$ws.Range("AJ70").Value = 'El usuario da clic en el botón "Imprimir".||El sistema muestra la descarga del archivo PDF de forma correcta.'

Set-CellStyle "A70" "A66"
Set-CellStyle "C70" "C2"
Set-CellStyle "D70" "D2"
Set-CellStyle "E70" "D2"
Set-CellStyle "F70" "F42"
Set-CellStyle "G70" "H38"
Set-CellStyle "H70" "H51"
Set-CellStyle "I70" "A2"
Set-CellStyle "J70" "A2"
Set-CellStyle "K70" "A2"
Set-CellStyle "L70" "A2"
Set-CellStyle "M70" "H10"
Set-CellStyle "N70" "H10"
Set-CellStyle "O70" "H10"
Set-CellStyle "P70" "H10"
Set-CellStyle "Q70" "H10"
Set-CellStyle "R70" "H10"
Set-CellStyle "S70" "H10"
Set-CellStyle "T70" "H10"
Set-CellStyle "U70" "H10"
Set-CellStyle "V70" "H10"
Set-CellStyle "W70" "U42"
Set-CellStyle "X70" "H10"
Set-CellStyle "Y70" "H10"
Set-CellStyle "Z70" "H10"
Set-CellStyle "AA70" "H10"
Set-CellStyle "AB70" "H10"
Set-CellStyle "AC70" "A2"
Set-CellStyle "AD70" "O15"
Set-CellStyle "AE70" "A2"
Set-CellStyle "AF70" "O15"
Set-CellStyle "AG70" "A2"
Set-CellStyle "AH70" "A2"
Set-CellStyle "AI70" "O15"
Set-CellStyle "AJ70" "O15"

# --- Row 71 ---
$ws.Range("A71").Value = 'TC_002_70_administradorConsulta'
$ws.Range("C71").Value = 'El usuario despliega la página de CCOP mediante el siguiente URL:
https://ccop-u.scointnet.net?countryId=MX||El sistema muestra la página de CCOP de forma correcta.'
$ws.Range("D71").Value = 'El usuario ingresa el siguiente dato: Nombre de usuario||El sistema permite ingresar el nombre de usuario de forma correcta.'
$ws.Range("E71").Value = 'El usuario ingresa el siguiente dato: Contraseña||El sistema permite ingresar la contraseña de forma correcta.'
$ws.Range("F71").Value = 'El usuario ingresa el siguiente dato: Token||El sistema permite ingresar el token de forma correcta.'
$ws.Range("G71").Value = 'El usuario da clic en el botón "Ingresar"||El sistema muestra la página principal de CCOP "Resumen consolidado de productos".'
$ws.Range("H71").Value = 'El usuario valida el total de cuentas mostradas entre paréntesis en la sección "LÍNEA OPERATIVA CONSOLIDADOS".||El sistema muestra el total de créditos de Línea Operativa de forma correcta.'
$ws.Range("I71").Value = 'El usuario da clic en la pestaña "Productos y servicios"||El sistema muestra la página "Cuentas".'
$ws.Range("J71").Value = 'El usuario da clic en el menú "Linea operativa"||El sistema muestra la lista de créditos de Línea Operativa y debe coincidir con el total de cuentas mostradas en el Dashboard en la sección "LÍNEA OPERATIVA CONSOLIDADOS".'
$ws.Range("K71").Value = 'El usuario valida que se muestre la opción "Ver consolidado"||El sistema muestra la opción "Ver consolidado" de forma correcta como hipervínculo.'
$ws.Range("L71").Value = 'El usuario da clic en la opción "Ver consolidado"||El sistema muestra la sección "LÍNEA OPERATIVA CONSOLIDADOS" con el número total de cuentas indicadas entre paréntesis y debe coincidir con el total de cuentas mostradas en el Dashboard en la sección "LÍNEA OPERATIVA CONSOLIDADOS".'
$ws.Range("M71").Value = 'El usuario valida el monto mostrado en el campo "Saldo total del crédito".||El sistema debe mostrar el mismo monto que el mostrado en el Dashboard en la sección "LÍNEA OPERATIVA CONSOLIDADOS"'
$ws.Range("N71").Value = 'El usuario valida el formato mostrado en el campo "Saldo total del crédito".||El sistema muestra el campo "Saldo total del crédito" de forma correcta con el siguiente formato:
Formato de moneda con signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("O71").Value = 'El usuario valida el monto mostrado en el campo "Saldo total del crédito".||El sistema debe mostrar el mismo monto que la suma de las cuentas de la columna "Saldo total del crédito".'
$ws.Range("P71").Value = 'El usuario valida el monto mostrado en el campo "Capital disponible".||El sistema debe mostrar el mismo monto que el mostrado en el Dashboard en la sección "LÍNEA OPERATIVA CONSOLIDADOS"'
$ws.Range("Q71").Value = 'El usuario valida el formato mostrado en el campo "Capital disponible".||El sistema muestra el campo "Capital disponible" de forma correcta con el siguiente formato:
Formato de moneda con signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("R71").Value = 'El usuario valida el monto mostrado en el campo "Capital disponible".||El sistema debe mostrar el mismo monto que la suma de las cuentas de la columna "Capital disponible".'
$ws.Range("S71").Value = 'El usuario valida el monto mostrado en el campo "Importe autorizado".||El sistema debe mostrar el mismo monto que el mostrado en el Dashboard en la sección "LÍNEA OPERATIVA CONSOLIDADOS"'
$ws.Range("T71").Value = 'El usuario valida el formato mostrado en el campo "Importe autorizado"||El sistema muestra el campo "Importe autorizado" de forma correcta con el siguiente formato:
Formato de moneda con signo $
Separación de miles y millones con separador de coma (,)
Signo de punto para decimales (.)
2 decimales'
$ws.Range("U71").Value = 'El usuario valida el monto mostrado en el campo "Importe autorizado".||El sistema debe mostrar el mismo monto que la suma de las cuentas de la columna "Importe autorizado".'

Set-CellStyle "A71" "A66"
Set-CellStyle "C71" "C2"
Set-CellStyle "D71" "D2"
Set-CellStyle "E71" "D2"
Set-CellStyle "F71" "F42"
Set-CellStyle "G71" "H38"
Set-CellStyle "H71" "A2"
Set-CellStyle "I71" "A2"
Set-CellStyle "J71" "A2"
Set-CellStyle "K71" "A2"
Set-CellStyle "L71" "A2"
Set-CellStyle "M71" "U42"
Set-CellStyle "N71" "H10"
Set-CellStyle "O71" "U42"
Set-CellStyle "P71" "H10"
Set-CellStyle "Q71" "H10"
Set-CellStyle "R71" "H10"
Set-CellStyle "S71" "H10"
Set-CellStyle "T71" "U42"
Set-CellStyle "U71" "H10"
